$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("C3").Value = "[-, 'MCT-3A-Automação Industrial', -, -]"
$ws.Range("D3").Value = "[-, -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F3").Value = "-"

# Row 4
$ws.Range("B4").Value = "[-, 'MCT-3A-Lab. Máquinas Elétricas', 'MCT-3A-Lab. Máquinas Elétricas']"
$ws.Range("C4").Value = "[-, 'MCT-3A-Automação Industrial', -, -]"
$ws.Range("D4").Value = "[-, -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F4").Value = "-"

# Row 6
$ws.Range("C6").Value = "[-, 'MCT-3A-Automação Industrial', -, -]"
$ws.Range("D6").Value = "[-, -, 'MCT-2A-Acionamentos Elétricos', -]"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("C7").Value = "[-, 'MCT-3A-Automação Industrial', -, -]"
$ws.Range("D7").Value = "-"
$ws.Range("F7").Value = "-"

# Row 8
$ws.Range("F8").Value = "['MCT-2A-Acionamentos Elétricos', -, -, -]"

# Row 18
$ws.Range("B18").Value = "[-, 'ELM-1NA-Manut. Elétrica', -, -]"
$ws.Range("C18").Value = "-"

# Row 19
$ws.Range("B19").Value = "[-, 'ELM-1NA-Manut. Elétrica', -, -]"
$ws.Range("C19").Value = "-"

# Row 20
$ws.Range("B20").Value = "[-, 'ELM-1NA-Manut. Elétrica', -, -]"
$ws.Range("C20").Value = "-"

# Row 21
$ws.Range("B21").Value = "[-, 'ELM-1NA-Manut. Elétrica', -, -]"
$ws.Range("C21").Value = "-"
